# Update the "Förändrad" (Changed) date column (C) for rows 2-27
# from serial date 45276 (2023-12-16) to 45277 (2023-12-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45276) {
        $cell.Value2 = 45277
    }
}
